$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interests")
$ws.Activate()

# New interest rows (id, interest) appended below the existing 14 rows
$newInterests = @(
    @(15, "drawing"),
    @(16, "sewing"),
    @(17, "cocktails"),
    @(18, "wine-tasting"),
    @(19, "genealogy"),
    @(20, "history"),
    @(21, "reading"),
    @(22, "calligraphy"),
    @(23, "critical thinking"),
    @(24, "poetry"),
    @(25, "youtube creation"),
    @(26, "web design"),
    @(27, "origami"),
    @(28, "entrepreneurship"),
    @(29, "geology"),
    @(30, "stamp collecting"),
    @(31, "candle making"),
    @(32, "jam making"),
    @(33, "knitting"),
    @(34, "goal-setting"),
    @(35, "gaming"),
    @(36, "gardening"),
    @(37, "sudoku"),
    @(38, "breathing exercises"),
    @(39, "comics"),
    @(40, "watercolours")
)

$startRow = 16
for ($i = 0; $i -lt $newInterests.Count; $i++) {
    $row = $startRow + $i
    $id = $newInterests[$i][0]
    $name = $newInterests[$i][1]
    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $name
}

# Update the active selection on the Interests sheet
$ws.Range("I6").Select()

# Update the workbook window position
$wb.Windows.Item(1).Left = 31950
$wb.Windows.Item(1).Top = 2760
